# naukri1.xlsx: append 4 more login rows (8-11) to Sheet1, mirroring the
# existing "email / password" rows, with a new email address and its own
# mailto hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$email = "marri24@gmail.com"
$password = "!1234567"

for ($r = 8; $r -le 11; $r++) {
    $emailCell = $ws.Cells.Item($r, 1)
    $passCell  = $ws.Cells.Item($r, 2)

    $emailCell.Value = $email
    $passCell.Value = $password

    $ws.Hyperlinks.Add($emailCell, "mailto:$email") | Out-Null

    $emailCell.Style = "Hyperlink"
    $passCell.Style = "Hyperlink"
}

$ws.Range("B11").Select() | Out-Null
